# Update "想去人数" (number of people interested) counts on the sheets
# that list 南宁·布谷鸟动漫展5th (row 3) and 南宁·2024良牙动漫秋季盛典 (row 4).
# These figures appear both on the "展览" sheet and the aggregated
# "全部类型" sheet, so both must be updated to stay consistent.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 329
    $ws.Range("F4").Value = 4547
}
